$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.252.19'
$ws.Cells.Item(2, 5).Value = '  +0.45%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.580.18'
$ws.Cells.Item(3, 5).Value = '  +1.70%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '315.60'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.97%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '97.07'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.88%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.55%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.05%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.53%  '

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '35.63'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.36%  '

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.0816'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.31%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '7.49'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -2.16%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '2.975.30'
$ws.Cells.Item(13, 5).Value = '  +1.67%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -3.62%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.581.03'
$ws.Cells.Item(15, 5).Value = '  +2.51%  '

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '15.17'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.57%  '

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.845'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.99%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '43.277.16'
$ws.Cells.Item(18, 5).Value = '  +0.48%  '

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +2.81%  '

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '12.61'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -3.84%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '0.0₃0964'
$ws.Cells.Item(21, 5).Value = '  -0.68%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '69.54'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.06%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '254.25'
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.61%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +2.30%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '27.15'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.09%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.01%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.22%  '

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '40.37'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.82%  '

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '10.32'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.04%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '5.85'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -3.85%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '154.76'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.07%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '3.41'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.71%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '2.14'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.13%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.0806'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.37%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.89%  '

# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '19.01'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.34%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.16%  '

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +6.04%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -0.76%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '22.53'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -5.84%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '3.97'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.58%  '

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.0304'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.81%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.05%  '

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '3.26'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -2.15%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.002.94'
$ws.Cells.Item(46, 5).Value = '  -1.22%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '8.93'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.90%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'BitcoinSV'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '83.31'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -3.32%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(49, 4).Value = '2.827.21'
$ws.Cells.Item(49, 5).Value = '  +1.69%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '74.94'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.95%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.194'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +1.80%  '
